$wb = $excel.ActiveWorkbook

# --- Rename "Sheet2" to "offlinechat" and populate it with the new chat/contact data ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "offlinechat"

# Header row
$ws2.Cells.Item(1,1).Value = "name"
$ws2.Cells.Item(1,2).Value = "email"
$ws2.Cells.Item(1,3).Value = "message"

# Data row
$ws2.Cells.Item(2,1).Value = "pavan"
$ws2.Cells.Item(2,2).Value = "pavandpagal@gmail.com"
$ws2.Cells.Item(2,3).Value = "hello how are you."

# Hyperlink the email address, then restore the canonical built-in "Hyperlink" cell
# style (Hyperlinks.Add otherwise mints its own style record).
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,2), "mailto:pavandpagal@gmail.com")
$ws2.Cells.Item(2,2).Style = "Hyperlink"

# Column widths for the new "email" / "message" columns
$ws2.Columns.Item(2).ColumnWidth = 22
$ws2.Columns.Item(3).ColumnWidth = 17.65

# This sheet becomes the active sheet/tab, with D3 selected
$ws2.Activate()
$ws2.Range("D3").Select() | Out-Null
